$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.285.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.06%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.931.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.11%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.13%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.7510"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.76%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'242.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.58%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.12%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'27.71"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.17%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.3174"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.93%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07099"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.05%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'TRON"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.08051"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.75%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'Polygon"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.7785"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.62%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.925.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.16%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.390"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.09%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'93.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.94%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'14.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'30.276.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'6.030"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +4.55%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'251.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.85%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.000007935"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.32%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'2.175.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.08%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.11%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.684"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.01%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'9.544"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.01%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.00%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'19.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.15%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.1298"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.27%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.189"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.57%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.367"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.98%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.551"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.62%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.412"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.44%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.143"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.36%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.05231"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.77%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.317"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +3.88%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.7583"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.77%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.781"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.33%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01953"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.38%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.797"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.05%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +2.46%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'78.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.13%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.4530"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.41%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.61%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.8406"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.82%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.10%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'9.989"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.35%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'7.669"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.20%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'101.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.03%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'38.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.76%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.1227"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +7.24%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'959.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.09%  "
$ws.Range("E51").Style = "Normal"
